# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The sheet's column G (header "K") holds per-row strike counts. This
# recomputes/rewrites those values (the stats such as std/mean that derive
# from them are implicitly refreshed since Excel recalculates on load).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for data rows 2..48 (row 1 is the header row).
$sVals = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 0
    14 = 1
    15 = 3
    16 = 2
    17 = 1
    18 = 1
    19 = 0
    20 = 2
    21 = 2
    22 = 0
    23 = 1
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 0
    36 = 1
    37 = 1
    38 = 0
    39 = 1
    40 = 0
    41 = 0
    42 = 1
    43 = 0
    44 = 1
    45 = 1
    46 = 0
    47 = 0
    48 = 1
}

foreach ($row in $sVals.Keys) {
    $ws.Range("G$row").Value = $sVals[$row]
}
